# Peppol Transport Profiles v8.4 - mark the legacy transport profiles as
# removed (TICC-248 / TICC-249) and stamp their removal date.
#
# Rows (Transport Profile sheet):
#   2 START  / busdox-transport-start
#   3 AS2 1.0 / busdox-transport-as2-ver1p0
#   4 AS4 1.0 / peppol-transport-as4-v1_0
#   5 AS4 2.0 / peppol-transport-as4-v2_0   (stays "active")
#   6 AS2 2.0 / busdox-transport-as2-ver2p0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transport Profile")

# AS2 2.0 (row 6) didn't carry a "Deprecation release" before - record it now.
$ws.Range("F6").Value() = "8.4"

# State: deprecated -> removed for every profile that is actually being
# removed in this release (row 5 / AS4 2.0 remains "active").
$ws.Range("E2").Value() = "removed"
$ws.Range("E3").Value() = "removed"
$ws.Range("E4").Value() = "removed"
$ws.Range("E6").Value() = "removed"

# Removal Date (column G) = 2023-08-24 for every removed profile,
# expressed as the Excel date serial number so no time-of-day fraction
# sneaks in.
$removalDate = 45162
$ws.Range("G2").Value() = $removalDate
$ws.Range("G3").Value() = $removalDate
$ws.Range("G4").Value() = $removalDate
$ws.Range("G6").Value() = $removalDate

# Reflect the author's final selection in the saved view.
$ws.Range("A7").Select()
